# Append a new blank data row (row 2) under the header row of the active
# sheet ("Đơn sale chính"). Only the numeric/financial columns get an
# explicit value (0); the remaining columns stay blank, matching the
# freshly-appended row produced by the report export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0   # Đơn giá gốc
$ws.Range("K2").Value = 0   # Upsale
$ws.Range("L2").Value = 0   # Đơn giá
$ws.Range("M2").Value = 0   # Thanh toán lần đầu
$ws.Range("N2").Value = 0   # Trả sau
$ws.Range("O2").Value = 0   # Đã thanh toán
$ws.Range("P2").Value = 0   # Dư nợ
